$wb = $excel.ActiveWorkbook

# The workbook currently has a single sheet "TC1". Add a new worksheet
# right after it named "ValidLogin" containing a small login-credentials
# table (username/password headers with admin/pointofsale values).
$tc1 = $wb.Worksheets.Item(1)
$validLogin = $wb.Worksheets.Add($null, $tc1)
$validLogin.Name = "ValidLogin"

$validLogin.Range("A1").Value = "username"
$validLogin.Range("B1").Value = "password"
$validLogin.Range("A2").Value = "admin"
$validLogin.Range("B2").Value = "pointofsale"

# Make the new sheet the active/visible one, zoomed to 160%, with B3
# selected as the next empty cell below the data.
[void]$validLogin.Select()
$excel.ActiveWindow.Zoom = 160
[void]$validLogin.Range("B3").Select()
